$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.066.60"
$ws.Range("E2").Value = "  +1.42%  "

# Row 3
$ws.Range("D3").Value = "3.331.00"

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'602.80"
$ws.Range("E5").Value = "  +1.35%  "

# Row 6
$ws.Range("D6").Value = "'144.21"
$ws.Range("E6").Value = "  +5.50%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "3.329.51"
$ws.Range("E8").Value = "  +6.55%  "

# Row 10
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = "  +3.29%  "

# Row 11
$ws.Range("D11").Value = "'5.57"
$ws.Range("E11").Value = "  +6.48%  "

# Row 12
$ws.Range("D12").Value = "'0.477"
$ws.Range("E12").Value = "  +4.13%  "

# Row 13
$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  +1.69%  "

# Row 14
$ws.Range("D14").Value = "'35.08"
$ws.Range("E14").Value = "  +2.34%  "

# Row 15
$ws.Range("D15").Value = "3.878.57"
$ws.Range("E15").Value = "  +6.37%  "

# Row 16
$ws.Range("E16").Value = "  +0.71%  "

# Row 17
$ws.Range("D17").Value = "3.323.82"
$ws.Range("E17").Value = "  +6.19%  "

# Row 18
$ws.Range("D18").Value = "64.158.65"
$ws.Range("E18").Value = "  +1.58%  "

# Row 19
$ws.Range("D19").Value = "'6.94"
$ws.Range("E19").Value = "  +3.48%  "

# Row 20
$ws.Range("D20").Value = "'485.01"
$ws.Range("E20").Value = "  +2.30%  "

# Row 21
$ws.Range("D21").Value = "'14.38"
$ws.Range("E21").Value = "  +0.85%  "

# Row 22
$ws.Range("D22").Value = "'0.743"
$ws.Range("E22").Value = "  +6.17%  "

# Row 23
$ws.Range("D23").Value = "'8.06"
$ws.Range("E23").Value = "  +4.16%  "

# Row 24
$ws.Range("D24").Value = "'13.86"
$ws.Range("E24").Value = "  +6.52%  "

# Row 25
$ws.Range("D25").Value = "'85.13"
$ws.Range("E25").Value = "  -1.74%  "

# Row 26
$ws.Range("E26").Value = "  +0.32%  "

# Row 27
$ws.Range("E27").Value = "  +2.36%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'8.38"
$ws.Range("E28").Value = "  +5.43%  "

# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'7.29"
$ws.Range("E29").Value = "  +1.92%  "

# Row 30
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.21%  "

# Row 31
$ws.Range("D31").Value = "'29.86"
$ws.Range("E31").Value = "  +11.20%  "

# Row 32
$ws.Range("E32").Value = "  +5.95%  "

# Row 33
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -1.21%  "

# Row 34
$ws.Range("D34").Value = "'2.59"
$ws.Range("E34").Value = "  +2.17%  "

# Row 35
$ws.Range("E35").Value = "  +2.61%  "

# Row 36
$ws.Range("D36").Value = "'6.04"
$ws.Range("E36").Value = "  +3.65%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0765"
$ws.Range("E37").Value = "  +7.85%  "

# Row 38
$ws.Range("D38").Value = "'53.43"
$ws.Range("E38").Value = "  +2.50%  "

# Row 39
$ws.Range("D39").Value = "'0.0405"
$ws.Range("E39").Value = "  +4.43%  "

# Row 40
$ws.Range("D40").Value = "'437.47"
$ws.Range("E40").Value = "  +3.39%  "

# Row 41
$ws.Range("D41").Value = "3.060.01"
$ws.Range("E41").Value = "  +5.71%  "

# Row 42
$ws.Range("D42").Value = "'2.82"
$ws.Range("E42").Value = "  +4.21%  "

# Row 43
$ws.Range("D43").Value = "'8.49"
$ws.Range("E43").Value = "  +2.94%  "

# Row 44
$ws.Range("D44").Value = "'0.112"
$ws.Range("E44").Value = "  -1.84%  "

# Row 45
$ws.Range("D45").Value = "'0.270"
$ws.Range("E45").Value = "  +2.94%  "

# Row 46
$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  +5.47%  "

# Row 47
$ws.Range("D47").Value = "'26.75"
$ws.Range("E47").Value = "  +3.90%  "

# Row 48
$ws.Range("D48").Value = "'36.27"
$ws.Range("E48").Value = "  +13.08%  "

# Row 50
$ws.Range("E50").Value = "  +2.73%  "

# Row 51
$ws.Range("E51").Value = "  +2.03%  "
